$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 20:12"
$ws.Range("B4").Value = 2915557
$ws.Range("C4").Value = 24969
$ws.Range("D4").Value = 1249330
$ws.Range("E4").Value = 1534018
$ws.Range("G4").Value = 145
$ws.Range("H4").Value = 132209
$ws.Range("B17").Value = 204610
$ws.Range("C17").Value = 1154
$ws.Range("D17").Value = 179492
$ws.Range("E17").Value = 19912
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 5206
$ws.Range("B65").Value = 13822
$ws.Range("C65").Value = 534
$ws.Range("D65").Value = 9329
$ws.Range("E65").Value = 4261
$ws.Range("B85").Value = 6159
$ws.Range("C85").Value = 101
$ws.Range("D85").Value = 4809
$ws.Range("E85").Value = 1298
$ws.Range("A91").Value = "Guayana Francesa"
$ws.Range("B91").Value = 4913
$ws.Range("C91").Value = 355
$ws.Range("D91").Value = 1866
$ws.Range("E91").Value = 3031
$ws.Range("H91").Value = 16
$ws.Range("A92").Value = "Republica de Yibuti"
$ws.Range("B92").Value = 4736
$ws.Range("D92").Value = 4580
$ws.Range("E92").Value = 101
$ws.Range("H92").Value = 55
$ws.Range("A93").Value = "Mauritania"
$ws.Range("B93").Value = 4705
$ws.Range("D93").Value = 1765
$ws.Range("E93").Value = 2811
$ws.Range("H93").Value = 129
$ws.Range("B103").Value = 2961
$ws.Range("C103").Value = 17
$ws.Range("D103").Value = 973
$ws.Range("E103").Value = 1896
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 92
$ws.Range("B108").Value = 2435
$ws.Range("C108").Value = 25
$ws.Range("E108").Value = 449
$ws.Range("B126").Value = 1421
$ws.Range("C126").Value = 39
$ws.Range("D126").Value = 654
$ws.Range("E126").Value = 751
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 16
$ws.Range("B128").Value = 1248
$ws.Range("C128").Value = 8
$ws.Range("D128").Value = 537
$ws.Range("E128").Value = 374
$ws.Range("G128").Value = 2
$ws.Range("H128").Value = 337
$ws.Range("A181").Value = "Aruba"
$ws.Range("B181").Value = 105
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 98
$ws.Range("H181").Value = 3
$ws.Range("A182").Value = "Bahamas"
$ws.Range("D182").Value = 89
$ws.Range("E182").Value = 4
$ws.Range("H182").Value = 11
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
